# Update the "Förändrad" (column C) date value on the active sheet.
# All data rows (2 through 81) had C = 45233 (2023-11-03); bump them
# forward to 45243 (2023-11-13), matching the upstream automatic
# file-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 81; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
